# Update countries & provincias Spain
# This script updates the COVID-19 data for a handful of countries and then
# re-sorts the data range by "Casos totales" (column B) descending, which is
# what produces the row reordering seen between Irlanda/Dinamarca and
# Egipto/Serbia/Hong Kong in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last used row in column A (data starts at row 4, header row 3)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Helper: find the row number whose column A text matches a given country name
function Find-CountryRow([string]$name) {
    $found = $ws.Range("A4:A$lastRow").Find($name, [Type]::Missing, [Type]::Missing, 1)
    if ($found -eq $null) {
        throw "Country '$name' not found"
    }
    return $found.Row
}

# New data values: Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes
$updates = @{
    "Estados Unidos" = @(98080, 12645, 2453, 94114, 2463, 218, 1513)
    "Alemania"       = @(50178, 6240, 6658, 43182, 23, 71, 338)
    "Irlanda"        = @(2121, 302, 5, 2094, 47, 3, 22)
    "Egipto"         = @(536, 41, 116, 390, 0, 6, 30)
    "Jordania"       = @(235, 23, 17, 218, 0, 0, 0)
}

foreach ($country in $updates.Keys) {
    $row = Find-CountryRow $country
    $vals = $updates[$country]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# Re-sort the data range (A4:H<lastRow>) by column B (Casos totales) descending,
# matching the site's ranking behaviour after the numbers were refreshed.
$sortRange = $ws.Range("A4:H$lastRow")
$key1 = $ws.Range("B4:B$lastRow")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 2, [Type]::Missing, 0) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142  # xlNo
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1  # xlTopToBottom
$ws.Sort.SortMethod = 1   # xlPinYin
$ws.Sort.Apply()
